$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.530.01"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.983.70"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "381.07"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").Value = "104.55"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("D7").Value = "0.548"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "37.30"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "3.454.84"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").Value = "18.45"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "2.988.71"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").Value = "0.974"
$ws.Range("E17").Value = "  +5.43%  "
$ws.Range("D18").Value = "51.507.80"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "7.46"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("D21").Value = "12.97"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "69.43"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").Value = "262.80"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "2.92"
$ws.Range("E25").Value = "  +9.78%  "
$ws.Range("D26").Value = "8.25"
$ws.Range("E26").Value = "  +16.50%  "
$ws.Range("D27").Value = "7.72"
$ws.Range("E27").Value = "  +23.08%  "
$ws.Range("E28").Value = "  +13.96%  "
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "9.90"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "34.92"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +7.60%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "17.20"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  +3.65%  "
$ws.Range("D43").Value = "125.47"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").Value = "21.76"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "0.281"
$ws.Range("E45").Value = "  +18.91%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("D49").Value = "2.033.57"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +7.71%  "
$ws.Range("D51").Value = "58.34"
$ws.Range("E51").Value = "  +3.81%  "
